# Weekly update: insert a new observation row at row 7 (pushing the
# existing rows 7-62 down to 8-63) and fill it with the latest data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7:62 down to 8:63, leaving row 7 blank (but keeping the D
# column's date style, since Insert copies formatting from the row above).
$ws.Rows("7:7").Insert()

# Populate the newly-inserted row 7 with the new weekly record.
$ws.Range("A7").Value = 9
$ws.Range("B7").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C7").Value = "Metropolitana"
$ws.Range("D7").Value = 44685
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = 100112029
$ws.Range("G7").Value = "Orégano"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 6
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 16000
$ws.Range("M7").Value = 15333
$ws.Range("N7").Value = '$/docena de atados'
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 5111
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = "Hortaliza"
